# This script reorders the roster rows (A2:C19) on the active sheet so that
# the player/position/team rows end up in the new order described by the
# commit's updated workbook, without altering the actual set of data or the
# header row/styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order, as (Player, Position, Team) tuples, for rows 2..19.
$data = @(
    @("Russell Westbrook", "PG",    "Denver Nuggets"),
    @("Jalen Green",       "PG,SG", "Houston Rockets"),
    @("Jalen Suggs",       "PG,SG", "Orlando Magic"),
    @("Chris Paul",        "PG",    "San Antonio Spurs"),
    @("Pascal Siakam",     "SF,PF", "Indiana Pacers"),
    @("Deni Avdija",       "SF,PF", "Portland Trail Blazers"),
    @("Naji Marshall",     "SG,SF", "Dallas Mavericks"),
    @("Jaylen Brown",      "SG,SF", "Boston Celtics"),
    @("Jakob Poeltl",      "C",     "Toronto Raptors"),
    @("Nikola Jokic",      "C",     "Denver Nuggets"),
    @("Dejounte Murray",   "PG,SG", "New Orleans Pelicans"),
    @("Jrue Holiday",      "PG,SG", "Boston Celtics"),
    @("Bogdan Bogdanovic", "SG,SF", "Atlanta Hawks"),
    @("Clint Capela",      "C",     "Atlanta Hawks"),
    @("Rudy Gobert",       "C",     "Minnesota Timberwolves"),
    @("Paolo Banchero",    "SF,PF", "Orlando Magic"),
    @("Chet Holmgren",     "PF,C",  "Oklahoma City Thunder"),
    @("Jerami Grant",      "SF,PF", "Portland Trail Blazers")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $triple = $data[$i]
    $ws.Cells.Item($row, 1).Value = $triple[0]
    $ws.Cells.Item($row, 2).Value = $triple[1]
    $ws.Cells.Item($row, 3).Value = $triple[2]
}
